$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: header + value (adds to sharedStrings + extends used range)
$ws.Range("D1").Value = "alerttext"
$ws.Range("D2").Value = "Customer added successfully"

# Column D width -> stored col width 25 (ColumnWidth uses Excel's padded
# character-width unit, which rounds to the stored width in the XML)
$ws.Columns.Item(4).ColumnWidth = 24.14

# Move the active selection as in the saved workbook
$ws.Range("E9").Select()
